# Adds "2022-Q3" quarter data:
#  1. Updates the "总计" (summary) sheet so a new top row for 2022-Q3 is
#     inserted and every later quarter row shifts down by one (2020-Q4
#     falls off the bottom of the old table onto a brand-new row).
#  2. Inserts a brand-new worksheet named "2022-Q3" (holding the per-fund
#     position breakdown for that quarter) right after "总计" and before
#     the former first quarter sheet ("2022-Q2").

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Update "总计" summary sheet
# ---------------------------------------------------------------------
$total = $wb.Worksheets.Item(1)

# Data rows, top to bottom, after the edit: (index, quarter, count, value)
$summaryRows = @(
    @(0, "2022-Q3", 6, 0.04),
    @(1, "2022-Q2", 2, 0.65),
    @(2, "2022-Q1", 1, 0.29),
    @(3, "2021-Q4", 1, 0.35),
    @(4, "2021-Q3", 2, 0.29),
    @(5, "2021-Q2", 6, 5.7),
    @(6, "2021-Q1", 3, 1.64),
    @(7, "2020-Q4", 1, 1.46)
)

# Carry the formatting of the previous last row (row 8) onto the newly
# used row 9 before filling in values.
$total.Range("A8").Copy()
$total.Range("A9").PasteSpecial(-4122)  # xlPasteFormats

for ($i = 0; $i -lt $summaryRows.Count; $i++) {
    $r = $i + 2
    $row = $summaryRows[$i]
    $total.Cells.Item($r, 1).Value = $row[0]
    $total.Cells.Item($r, 2).Value = $row[1]
    $total.Cells.Item($r, 3).Value = $row[2]
    $total.Cells.Item($r, 4).Value = $row[3]
}

# ---------------------------------------------------------------------
# 2. Insert the new "2022-Q3" worksheet
# ---------------------------------------------------------------------
$q2 = $wb.Worksheets.Item("2022-Q2")
$q3 = $wb.Worksheets.Add($q2, $null)
$q3.Name = "2022-Q3"

# Copy header/data-row formatting from the "2022-Q2" sheet so the new
# sheet matches the look of its siblings (bold bordered header row and
# bold bordered index column).
$q2.Range("B1:H1").Copy()
$q3.Range("B1:H1").PasteSpecial(-4122)  # xlPasteFormats
$q2.Range("A2").Copy()
$q3.Range("A2:A7").PasteSpecial(-4122)  # xlPasteFormats

$headers = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
for ($c = 0; $c -lt $headers.Count; $c++) {
    $q3.Cells.Item(1, $c + 2).Value = $headers[$c]
}

# (index, code, name, scale, position, ratio, marketValue, rank)
$fundRows = @(
    @(0, "620002", "金元顺安成长动力混合", "0.51", "62.21", "2.96", "0.0151", 9),
    @(1, "005381", "泰康睿利量化多策略混合A", "0.48", "78.29", "2.03", "0.0097", 8),
    @(2, "005382", "泰康睿利量化多策略混合C", "0.48", "78.29", "2.03", "0.0097", 8),
    @(3, "009327", "东兴兴晟混合A", "0.38", "79.70", "1.14", "0.0043", 4),
    @(4, "002952", "建信多因子量化股票", "0.09", "91.26", "3.23", "0.0029", 7),
    @(5, "009328", "东兴兴晟混合C", "0.07", "79.70", "1.14", "0.0008", 4)
)

for ($i = 0; $i -lt $fundRows.Count; $i++) {
    $r = $i + 2
    $row = $fundRows[$i]
    $q3.Cells.Item($r, 1).Value = $row[0]
    $q3.Cells.Item($r, 2).Value = "'" + $row[1]
    $q3.Cells.Item($r, 3).Value = $row[2]
    $q3.Cells.Item($r, 4).Value = "'" + $row[3]
    $q3.Cells.Item($r, 5).Value = "'" + $row[4]
    $q3.Cells.Item($r, 6).Value = "'" + $row[5]
    $q3.Cells.Item($r, 7).Value = "'" + $row[6]
    $q3.Cells.Item($r, 8).Value = $row[7]
}
